# "Tweaks after user testing."
#  - bump the cached date field text (master + every layout) by a few days
#  - re-brand the purple "Tip:" font from FuturaHandwritten -> Articulate
#  - re-brand the title font from Peace Sans -> Articulate Extrabold
#  - make sure every "Tip: " lead-in is bold + purple + Articulate Extra... font,
#    splitting it out of the body run where it was still merged in

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Cached "datetimeFigureOut" field text: 4/28/2020 -> 5/2/2020
#    (lives on the slide master's Date placeholder and on each of the
#    11 layouts' own copy of it)
# ---------------------------------------------------------------------------
$oldDate = "4/28/2020"
$newDate = "5/2/2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Title on slide 1: Peace Sans -> Articulate Extrabold
# ---------------------------------------------------------------------------
$titleRange = $p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange
$titleRange.Font.Name = "Articulate Extrabold"

# ---------------------------------------------------------------------------
# 3) Every "Tip: ..." textbox (TextBox 1, shape 5) on slides 2-7:
#       - "Tip: " lead-in -> bold, purple (7030A0), Articulate
#       - remainder of the text -> Articulate (keep existing formatting)
#    On slides 4 and 7 the lead-in was still merged into the body run, so
#    splitting by character range also takes care of separating it out.
# ---------------------------------------------------------------------------
$purple = 0xA03070  # BGR-packed COM value for srgbClr 7030A0

for ($slideIdx = 2; $slideIdx -le 7; $slideIdx++) {
    $tipBox = $p.Slides.Item($slideIdx).Shapes.Item(5)
    $tr = $tipBox.TextFrame.TextRange
    $fullText = $tr.Text
    $leadLen = 5  # Length of "Tip: "

    $lead = $tr.Characters(1, $leadLen)
    $lead.Font.Bold = $true
    $lead.Font.Color.RGB = $purple
    $lead.Font.Name = "Articulate"

    $remainingLen = $fullText.Length - $leadLen
    if ($remainingLen -gt 0) {
        $body = $tr.Characters($leadLen + 1, $remainingLen)
        $body.Font.Name = "Articulate"
    }
}
